# Update "paises.xlsx" (sheet "Pais") with the latest COVID-19 country
# statistics snapshot and refresh the "last updated" timestamp.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Row 1: "Datos actualizados a ..." timestamp banner -> bump time 03:32 -> 04:49
$ws.Range("A1").Value = "Datos actualizados a 24 de Septiembre de 2020 a las 04:49"

# Row 50: Honduras (Casos totales, Nuevos casos, Casos activos, Recuperados, Casos criticos, Muertes hoy, Muertes)
$ws.Range("B50").Value = 72675
$ws.Range("C50").Value = 369
$ws.Range("D50").Value = 24022
$ws.Range("E50").Value = 46431
$ws.Range("G50").Value = 16
$ws.Range("H50").Value = 2222

# Row 53: Venezuela
$ws.Range("B53").Value = 69439
$ws.Range("D53").Value = 58759
$ws.Range("E53").Value = 10106
$ws.Range("H53").Value = 574

# Row 72: Paraguay
$ws.Range("B72").Value = 35571
$ws.Range("D72").Value = 19867
$ws.Range("E72").Value = 14977
$ws.Range("H72").Value = 727

# Row 78: Australia
$ws.Range("B78").Value = 26980
$ws.Range("C78").Value = 6
$ws.Range("D78").Value = 24417
$ws.Range("E78").Value = 1702
$ws.Range("G78").Value = 2
$ws.Range("H78").Value = 861

# Row 108: Zimbabue
$ws.Range("B108").Value = 7827
$ws.Range("C108").Value = 535
$ws.Range("D108").Value = 2085
$ws.Range("E108").Value = 5609
$ws.Range("G108").Value = 3
$ws.Range("H108").Value = 133

# Row 109: Mauritania
$ws.Range("B109").Value = 7725
$ws.Range("D109").Value = 6007
$ws.Range("E109").Value = 1491
$ws.Range("H109").Value = 227

# Row 110: Eslovaquia
$ws.Range("B110").Value = 7425
$ws.Range("D110").Value = 7028
$ws.Range("E110").Value = 236
$ws.Range("H110").Value = 161

# Row 111: Mozambique
$ws.Range("B111").Value = 7269
$ws.Range("D111").Value = 3888
$ws.Range("E111").Value = 3340
$ws.Range("H111").Value = 41

# Row 112: Birmania
$ws.Range("B112").Value = 7262
$ws.Range("D112").Value = 4350
$ws.Range("E112").Value = 2863
$ws.Range("H112").Value = 49

# Row 156: Nueva Zelanda
$ws.Range("B156").Value = 1827
$ws.Range("C156").Value = 3
$ws.Range("E156").Value = 65

# Row 163: Polinesia Francesa
$ws.Range("B163").Value = 1469
$ws.Range("C163").Value = 75
$ws.Range("D163").Value = 1237
$ws.Range("E163").Value = 227
$ws.Range("G163").Value = 3
$ws.Range("H163").Value = 5

# Row 214: Montserrat
$ws.Range("D214").Value = 13
$ws.Range("H214").Value = 0

# Row 215: Islas Malvinas
$ws.Range("D215").Value = 12
$ws.Range("H215").Value = 1
